$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header G1: pop_sq_mile_1mi -> pop_sq_mile_3mi
$ws.Range("G1").Value = "pop_sq_mile_3mi"

# Column R and S value updates
$ws.Range("R2").Value = 49
$ws.Range("S3").Value = 0.5
$ws.Range("R4").Value = 34.5454545454545
$ws.Range("S4").Value = 0.3
$ws.Range("S5").Value = 0.373913043478261
$ws.Range("R6").Value = 28.974358974359
$ws.Range("S6").Value = 0.361538461538462
$ws.Range("R7").Value = 30
$ws.Range("S7").Value = 0.342857142857143
$ws.Range("R8").Value = 62.5
$ws.Range("S8").Value = 0.4375
$ws.Range("R9").Value = 61
$ws.Range("S9").Value = 0.47
$ws.Range("R10").Value = 20
$ws.Range("S10").Value = 0.2
$ws.Range("R11").Value = 19
